$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Fitness values (column C) for rows 2-46
$values = @(13543,12795,12762,12330,12330,11958,11958,10208,10036,10036,10036,10036,10036,10036,9883,9493,9045,9045,9020,8971,8971,8780,8780,8586,8296,8296,8296,8279,7954,7954,7929,7929,7929,7929,7929,7929,7929,7929,7929,7929,7929,7929,7678,7678,7678)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
